# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  "Office Theme" (Office colours)  -> used by the Notes Master
#   ppt/theme/theme2.xml  "Integral"     (Red Violet)       -> used by the Slide Master / Design
#
# The target revision swaps the *content* of the two theme parts (the
# filenames / relationships are untouched): theme2.xml (the Slide Master's
# theme, reachable through the PowerPoint object model) ends up holding the
# plain "Office Theme" colour palette that used to live in theme1.xml.
#
# PowerPoint's COM object model doesn't provide a raw "set this OOXML part's
# bytes" call, so we reproduce the effect through ThemeColorScheme, which is
# the supported, documented way to repaint a design's 12 theme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that order).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (the colours that theme1.xml originally had),
# expressed as OLE (BGR) integers for the RGB property setter.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
